$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two existing hyperlinks that cover D3:D21 and D22:D27 - they will
# be replaced by a single hyperlink after the new row is inserted.
$ws.Hyperlinks.Item(3).Delete()
$ws.Hyperlinks.Item(2).Delete()

# Insert a new row at row 15 (shifts existing rows 15-27 down to 16-28),
# pushing "United Kingdom" in ahead of the rest-of-world countries list.
$ws.Rows.Item(15).Insert()

# Fill in the new row's data: date, country, cases, source URL text
$ws.Cells.Item(15, 1).Value = 44691
$ws.Cells.Item(15, 2).Value = "United Kingdom"
$ws.Cells.Item(15, 3).Value = 163
$ws.Cells.Item(15, 4).Value = "https://www.ecdc.europa.eu/en/news-events/epidemiological-update-hepatitis-unknown-aetiology-children"

# Copy formatting (number format / cell style) from row 14 into the new row
# 15 for the date and source columns (A and D).
$ws.Cells.Item(14, 1).Copy() | Out-Null
$ws.Cells.Item(15, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14, 4).Copy() | Out-Null
$ws.Cells.Item(15, 4).PasteSpecial(-4122) | Out-Null

# Re-create the combined hyperlink over the now-shifted D23:D28 range.
$ws.Hyperlinks.Add($ws.Range("D23:D28"), "https://www.ecdc.europa.eu/en/news-events/epidemiological-update-hepatitis-unknown-aetiology-children", "", "", "https://www.ecdc.europa.eu/en/news-events/epidemiological-update-hepatitis-unknown-aetiology-children") | Out-Null

# Adding the hyperlink re-stamps the anchor cell (D23) with a fresh style
# variant; normalize it back to the same style used by the rest of the
# (unchanged) Source column cells.
$ws.Cells.Item(24, 4).Copy() | Out-Null
$ws.Cells.Item(23, 4).PasteSpecial(-4122) | Out-Null

# Widen column B to fit the new "United Kingdom" entry (closest attainable
# width to Excel's computed best-fit for this content).
$ws.Columns.Item(2).ColumnWidth = 11

# Extend the table range / autofilter to include the new row.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:D28")) | Out-Null

$ws.Range("B16").Select()
